$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 4281
$ws.Range("I62").Value = 2452.5
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 2452.5
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -1828.5
$ws.Range("N62").Value = -6748
$ws.Range("H65").Value = 4281
$ws.Range("I65").Value = 2452.5
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 12262.5
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -9142.5
$ws.Range("N65").Value = -33740
$ws.Range("H86").Value = 4932.885
$ws.Range("I86").Value = 7333.6665
$ws.Range("J86").Value = 4619.7393
$ws.Range("K86").Value = 7333.6665
$ws.Range("L86").Value = 4619.7393
$ws.Range("M86").Value = -6210.6665
$ws.Range("N86").Value = -6865.7393
$ws.Range("H89").Value = 4932.885
$ws.Range("I89").Value = 7333.6665
$ws.Range("J89").Value = 4619.7393
$ws.Range("K89").Value = 36668.3325
$ws.Range("L89").Value = 23098.6965
$ws.Range("M89").Value = -31052.3325
$ws.Range("N89").Value = -34330.69650000001
$ws.Range("H111").Value = 2916
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 2916
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 8748
$ws.Range("M111").ClearContents()
$ws.Range("N111").Value = -14882
$ws.Range("H129").Value = 965.56757
$ws.Range("I129").Value = 375
$ws.Range("J129").Value = 1037.1515
$ws.Range("K129").Value = 1125
$ws.Range("L129").Value = 3111.4545
$ws.Range("M129").Value = 3875
$ws.Range("N129").Value = -13111.4545

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6228.1357
$ws.Range("I32").Value = 5241.269
$ws.Range("J32").Value = 13559.143
$ws.Range("K32").Value = 5241.269
$ws.Range("L32").Value = 13559.143
$ws.Range("M32").Value = -4954.269
$ws.Range("N32").Value = -14133.143
$ws.Range("H45").Value = 1459.8572
$ws.Range("I45").Value = 1478.7693
$ws.Range("J45").Value = 1214
$ws.Range("K45").Value = 1478.7693
$ws.Range("L45").Value = 1214
$ws.Range("M45").Value = -1101.7693
$ws.Range("N45").Value = -1968
$ws.Range("H97").Value = 1062.05
$ws.Range("I97").Value = 302.94116
$ws.Range("J97").Value = 5363.6665
$ws.Range("K97").Value = 302.94116
$ws.Range("L97").Value = 5363.6665
$ws.Range("M97").Value = 193.05884
$ws.Range("N97").Value = -6355.6665
$ws.Range("H122").Value = 1537.5186
$ws.Range("I122").Value = 1109.7368
$ws.Range("J122").Value = 2553.5
$ws.Range("K122").Value = 3329.2104
$ws.Range("L122").Value = 7660.5
$ws.Range("M122").Value = -879.2103999999999
$ws.Range("N122").Value = -12560.5
$ws.Range("H132").Value = 6854.625
$ws.Range("I132").Value = 10401
$ws.Range("J132").Value = 2663.4546
$ws.Range("K132").Value = 31203
$ws.Range("L132").Value = 7990.3638
$ws.Range("M132").Value = -28673
$ws.Range("N132").Value = -13050.3638

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1457.3715
$ws.Range("I20").Value = 1509.625
$ws.Range("J20").Value = 900
$ws.Range("K20").Value = 1509.625
$ws.Range("L20").Value = 900
$ws.Range("M20").Value = -1262.625
$ws.Range("N20").Value = -1394
$ws.Range("H107").Value = 707.4545000000001
$ws.Range("I107").Value = 529.8889
$ws.Range("K107").Value = 529.8889
$ws.Range("M107").Value = 1390.1111
$ws.Range("H134").Value = 4449.4116
$ws.Range("I134").Value = 5037.037
$ws.Range("J134").Value = 2182.8572
$ws.Range("K134").Value = 15111.111
$ws.Range("L134").Value = 6548.571599999999
$ws.Range("M134").Value = -12576.111
$ws.Range("N134").Value = -11618.5716

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33511.125
$ws.Range("I31").Value = 60270.94
$ws.Range("K31").Value = 60270.94
$ws.Range("M31").Value = -59975.94
$ws.Range("H34").Value = 33511.125
$ws.Range("I34").Value = 60270.94
$ws.Range("K34").Value = 60270.94
$ws.Range("M34").Value = -60068.94
$ws.Range("H58").Value = 1165.4166
$ws.Range("I58").Value = 1178
$ws.Range("J58").Value = 1110.8889
$ws.Range("K58").Value = 1178
$ws.Range("L58").Value = 1110.8889
$ws.Range("M58").Value = -975
$ws.Range("N58").Value = -1516.8889
$ws.Range("H99").Value = 1489.4706
$ws.Range("I99").Value = 1507.625
$ws.Range("J99").Value = 1473.3334
$ws.Range("K99").Value = 1507.625
$ws.Range("L99").Value = 1473.3334
$ws.Range("M99").Value = -9.625
$ws.Range("N99").Value = -4469.3334
$ws.Range("H126").Value = 1489.4706
$ws.Range("I126").Value = 1507.625
$ws.Range("J126").Value = 1473.3334
$ws.Range("K126").Value = 4522.875
$ws.Range("L126").Value = 4420.0002
$ws.Range("M126").Value = -2052.875
$ws.Range("N126").Value = -9360.0002
$ws.Range("H136").Value = 1165.4166
$ws.Range("I136").Value = 1178
$ws.Range("J136").Value = 1110.8889
$ws.Range("K136").Value = 3534
$ws.Range("L136").Value = 3332.6667
$ws.Range("M136").Value = -984
$ws.Range("N136").Value = -8432.6667

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1046.8
$ws.Range("J5").Value = 1167.8948
$ws.Range("L5").Value = 3503.6844
$ws.Range("N5").Value = -3727.6844
$ws.Range("H116").Value = 5739.231
$ws.Range("I116").Value = 909.3333
$ws.Range("J116").Value = 7188.2
$ws.Range("K116").Value = 2727.9999
$ws.Range("L116").Value = 21564.6
$ws.Range("M116").Value = 714.0001000000002
$ws.Range("N116").Value = -28448.6
$ws.Range("H118").Value = 2231.0454
$ws.Range("I118").Value = 671
$ws.Range("J118").Value = 2689.8823
$ws.Range("K118").Value = 2013
$ws.Range("L118").Value = 8069.646900000001
$ws.Range("M118").Value = -770
$ws.Range("N118").Value = -10555.6469
$ws.Range("H135").Value = 1046.8
$ws.Range("J135").Value = 1167.8948
$ws.Range("L135").Value = 10511.0532
$ws.Range("N135").Value = -15581.0532

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4946.6665
$ws.Range("I70").Value = 4741.75
$ws.Range("J70").Value = 5202.8125
$ws.Range("K70").Value = 4741.75
$ws.Range("L70").Value = 5202.8125
$ws.Range("M70").Value = -4471.75
$ws.Range("N70").Value = -5742.8125
$ws.Range("H73").Value = 4946.6665
$ws.Range("I73").Value = 4741.75
$ws.Range("J73").Value = 5202.8125
$ws.Range("K73").Value = 4741.75
$ws.Range("L73").Value = 5202.8125
$ws.Range("M73").Value = -3805.75
$ws.Range("N73").Value = -7074.8125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 7200.5557
$ws.Range("I122").Value = 7828.5713
$ws.Range("J122").Value = 5002.5
$ws.Range("K122").Value = 23485.7139
$ws.Range("L122").Value = 15007.5
$ws.Range("M122").Value = -21035.7139
$ws.Range("N122").Value = -19907.5
$ws.Range("H132").Value = 3207.1345
$ws.Range("I132").Value = 2693.1143
$ws.Range("J132").Value = 4265.4116
$ws.Range("K132").Value = 8079.342900000001
$ws.Range("L132").Value = 12796.2348
$ws.Range("M132").Value = -5549.342900000001
$ws.Range("N132").Value = -17856.2348
$ws.Range("H136").Value = 2209.1353
$ws.Range("I136").Value = 1740.2554
$ws.Range("J136").Value = 3025.3333
$ws.Range("K136").Value = 5220.7662
$ws.Range("L136").Value = 9075.999899999999
$ws.Range("M136").Value = -2670.7662
$ws.Range("N136").Value = -14175.9999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1503777
$ws.Range("I126").Value = 1161338.1
$ws.Range("J126").Value = 2416947.2
$ws.Range("K126").Value = 3484014.3
$ws.Range("L126").Value = 7250841.600000001
$ws.Range("M126").Value = -3481544.3
$ws.Range("N126").Value = -7255781.600000001
